$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.860.56'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.552.47'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.35'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.47'
$ws.Range("E8").Value = '  -3.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0583'
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0856'
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.773.11'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.562.04'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.853.86'
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.34'
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.57'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("D19").Value = '0.0₃0689'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.10'
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.01'
$ws.Range("E23").Value = '  -4.09%  '
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.55'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.54'
$ws.Range("E26").Value = '  -2.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.95'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.354.60'
$ws.Range("E33").Value = '  -3.94%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.52'
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.934'
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.66'
$ws.Range("E42").Value = '  +5.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.993'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.29'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.686.29'
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("E48").Value = '  -3.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.02'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0508'
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("E51").Value = '  -2.49%  '
